# edit.ps1
# Applies the cryptos list update described by the commit:
# "Updated cryptos list on Wed Sep 27 19:30:39 UTC 2023 with GitHub Actions"
#
# The sheet holds a list of cryptocurrencies with columns:
#   A = rank index, B = Coin name, C = Link, D = Price, E = Volume(1h)
# This update refreshes the Price/Volume figures for most rows, and for
# rows 50-51 also inserts a new coin (Algorand) ahead of USDD, shifting
# USDD down into row 51 (replacing the former EnergySwap row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "26.269.67" },
    @{ Cell = "E2"; Value = "  +0.28%  " },
    @{ Cell = "D3"; Value = "1.597.23" },
    @{ Cell = "E3"; Value = "  +0.72%  " },
    @{ Cell = "E4"; Value = "  -0.05%  " },
    @{ Cell = "D5"; Value = "211.85" },
    @{ Cell = "E5"; Value = "  +0.04%  " },
    @{ Cell = "E6"; Value = "  -0.01%  " },
    @{ Cell = "E7"; Value = "  -0.04%  " },
    @{ Cell = "D8"; Value = "0.245" },
    @{ Cell = "E8"; Value = "  +0.22%  " },
    @{ Cell = "D9"; Value = "0.0606" },
    @{ Cell = "E9"; Value = "  +0.43%  " },
    @{ Cell = "D10"; Value = "19.01" },
    @{ Cell = "E10"; Value = "  -0.87%  " },
    @{ Cell = "E11"; Value = "  +0.75%  " },
    @{ Cell = "D12"; Value = "1.820.33" },
    @{ Cell = "E12"; Value = "  +0.64%  " },
    @{ Cell = "D13"; Value = "1.601.76" },
    @{ Cell = "E13"; Value = "  +0.76%  " },
    @{ Cell = "E14"; Value = "  -0.07%  " },
    @{ Cell = "D15"; Value = "0.504" },
    @{ Cell = "E15"; Value = "  -2.29%  " },
    @{ Cell = "D16"; Value = "63.73" },
    @{ Cell = "E16"; Value = "  -0.21%  " },
    @{ Cell = "D17"; Value = "26.252.82" },
    @{ Cell = "E17"; Value = "  +0.26%  " },
    @{ Cell = "D18"; Value = "230.42" },
    @{ Cell = "E18"; Value = "  +7.75%  " },
    @{ Cell = "E19"; Value = "  +5.20%  " },
    @{ Cell = "D20"; Value = "0.0₃0723" },
    @{ Cell = "E20"; Value = "  -0.19%  " },
    @{ Cell = "D21"; Value = "1.00" },
    @{ Cell = "E21"; Value = "  +0.04%  " },
    @{ Cell = "D22"; Value = "4.24" },
    @{ Cell = "E22"; Value = "  -0.09%  " },
    @{ Cell = "E23"; Value = "  +1.42%  " },
    @{ Cell = "D24"; Value = "8.94" },
    @{ Cell = "E24"; Value = "  -0.36%  " },
    @{ Cell = "D25"; Value = "145.84" },
    @{ Cell = "E25"; Value = "  +0.95%  " },
    @{ Cell = "E26"; Value = "  -0.01%  " },
    @{ Cell = "D27"; Value = "7.03" },
    @{ Cell = "E27"; Value = "  +0.78%  " },
    @{ Cell = "E28"; Value = "  +0.54%  " },
    @{ Cell = "D29"; Value = "15.34" },
    @{ Cell = "E29"; Value = "  +1.72%  " },
    @{ Cell = "D30"; Value = "0.0493" },
    @{ Cell = "E30"; Value = "  -0.38%  " },
    @{ Cell = "E31"; Value = "  +0.36%  " },
    @{ Cell = "D32"; Value = "3.21" },
    @{ Cell = "E32"; Value = "  +0.97%  " },
    @{ Cell = "D33"; Value = "1.466.61" },
    @{ Cell = "E33"; Value = "  +3.83%  " },
    @{ Cell = "E34"; Value = "  +0.42%  " },
    @{ Cell = "E35"; Value = "  -0.32%  " },
    @{ Cell = "D36"; Value = "1.47" },
    @{ Cell = "E36"; Value = "  +0.93%  " },
    @{ Cell = "D37"; Value = "0.570" },
    @{ Cell = "E37"; Value = "  -2.85%  " },
    @{ Cell = "E38"; Value = "  -0.94%  " },
    @{ Cell = "D39"; Value = "0.822" },
    @{ Cell = "E39"; Value = "  +0.25%  " },
    @{ Cell = "E40"; Value = "  -2.04%  " },
    @{ Cell = "E41"; Value = "  +0.02%  " },
    @{ Cell = "E42"; Value = "  +2.09%  " },
    @{ Cell = "D43"; Value = "0.932" },
    @{ Cell = "E43"; Value = "  -2.07%  " },
    @{ Cell = "D44"; Value = "1.732.45" },
    @{ Cell = "E44"; Value = "  +0.71%  " },
    @{ Cell = "D45"; Value = "0.758" },
    @{ Cell = "E45"; Value = "  -0.92%  " },
    @{ Cell = "D46"; Value = "60.61" },
    @{ Cell = "E46"; Value = "  -0.53%  " },
    @{ Cell = "D47"; Value = "87.71" },
    @{ Cell = "E47"; Value = "  +2.84%  " },
    @{ Cell = "E48"; Value = "  -0.13%  " },
    @{ Cell = "E49"; Value = "  +0.10%  " },
    @{ Cell = "B50"; Value = "Algorand" },
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" },
    @{ Cell = "D50"; Value = "0.0948" },
    @{ Cell = "E50"; Value = "  -2.03%  " },
    @{ Cell = "B51"; Value = "USDD" },
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd" },
    @{ Cell = "D51"; Value = "0.997" },
    @{ Cell = "E51"; Value = "  -0.10%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # The Price column (D) sometimes holds plain decimal-looking strings
    # (e.g. "211.85", "1.00", "0.997"). Excel's automatic type detection
    # would otherwise silently convert these to numbers (losing trailing
    # zeros / exact text), so force the cell to Text format first, same
    # as the rest of the sheet's inline-string cells.
    if ($u.Cell -match '^D\d+$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u.Value
}
